$wb = $excel.ActiveWorkbook

# Regression-testing run timestamps: each of these sheets has an
# execution-log column B (rows 2-5) that records when the row was last
# run. This run refreshes the stamps to the latest regression pass.

$ws = $wb.Worksheets.Item("PayNowCC")
$ws.Range("B2").Value = "Fri Nov 14 22:00:51 IST 2025"
$ws.Range("B3").Value = "Fri Nov 14 22:02:04 IST 2025"
$ws.Range("B4").Value = "Fri Nov 14 22:03:15 IST 2025"
$ws.Range("B5").Value = "Fri Nov 14 22:04:45 IST 2025"

$ws = $wb.Worksheets.Item("PayNowCCSCF")
$ws.Range("B2").Value = "Fri Nov 14 22:05:44 IST 2025"
$ws.Range("B3").Value = "Fri Nov 14 22:06:46 IST 2025"
$ws.Range("B4").Value = "Fri Nov 14 22:07:52 IST 2025"
$ws.Range("B5").Value = "Fri Nov 14 22:08:56 IST 2025"

$ws = $wb.Worksheets.Item("PayNowCCDCF")
$ws.Range("B2").Value = "Fri Nov 14 22:09:59 IST 2025"
$ws.Range("B3").Value = "Fri Nov 14 22:11:01 IST 2025"
$ws.Range("B4").Value = "Fri Nov 14 22:12:01 IST 2025"
$ws.Range("B5").Value = "Fri Nov 14 22:13:05 IST 2025"

$ws = $wb.Worksheets.Item("NoModifyAmount")
$ws.Range("B2").Value = "Fri Nov 14 22:15:55 IST 2025"

$ws = $wb.Worksheets.Item("OverUnderPay")
$ws.Range("B2").Value = "Fri Nov 14 22:22:27 IST 2025"
$ws.Range("B3").Value = "Fri Nov 14 22:23:04 IST 2025"

$ws = $wb.Worksheets.Item("NoOverPay")
$ws.Range("B2").Value = "Fri Nov 14 22:26:45 IST 2025"
